# Weekly update: a new price-survey entry is published for the
# "Feria Lagunitas de Puerto Montt - Betarraga" sheet. The new record is
# inserted as the first data row of the date-ordered block (row 284),
# pushing every existing row in that block down by one (284->285, ...,
# 363->364). This is the same pattern Excel's own Insert-a-row does, so we
# simply insert a blank row at 284 (which shifts 284:363 down to 285:364,
# carrying their original values/styles with them) and then populate the
# freshly inserted row with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 284; rows 284-363 shift down to 285-364.
$ws.Rows("284:284").Insert()

# Populate the newly inserted row 284 with the new weekly record.
$ws.Cells.Item(284, 1).Value = 4
$ws.Cells.Item(284, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(284, 3).Value = "Los Lagos"
$ws.Cells.Item(284, 4).Value = 44841
$ws.Cells.Item(284, 5).Value = 10
$ws.Cells.Item(284, 6).Value = 100114014
$ws.Cells.Item(284, 7).Value = "Betarraga"
$ws.Cells.Item(284, 8).Value = "Sin especificar"
$ws.Cells.Item(284, 9).Value = "Primera"
$ws.Cells.Item(284, 10).Value = 1000
$ws.Cells.Item(284, 11).Value = 1500
$ws.Cells.Item(284, 12).Value = 1500
$ws.Cells.Item(284, 13).Value = 1500
$ws.Cells.Item(284, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(284, 15).Value = "Región del Maule"
$ws.Cells.Item(284, 16).Value = 300
$ws.Cells.Item(284, 17).Value = 5
$ws.Cells.Item(284, 18).Value = "Hortaliza"
